# Updates the "Hot Stock Top 20" table (columns A:C, rows 2:21) with the
# refreshed rankings for 财联社 (A), 东方财富 (B) and 同花顺 (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("蓝色光标", "蓝色光标", "蓝色光标"),
    @("利欧股份", "利欧股份", "利欧股份"),
    @("岩山科技", "岩山科技", "航天发展"),
    @("航天发展", "中国卫星", "天际股份"),
    @("中国卫星", "航天发展", "岩山科技"),
    @("创新医疗", "金风科技", "航天机电"),
    @("金风科技", "中国卫通", "雷科防务"),
    @("中国卫通", "雷科防务", "金风科技"),
    @("三花智控", "三花智控", "海格通信"),
    @("雷科防务", "创新医疗", "中国卫通"),
    @("三博脑科", "航天电子", "中国卫星"),
    @("南京熊猫", "天际股份", "平潭发展"),
    @("天际股份", "三博脑科", "泰尔股份"),
    @("航天机电", "航天机电", "神剑股份"),
    @("海格通信", "南京熊猫", "三花智控"),
    @("北斗星通", "北斗星通", "国晟科技"),
    @("烽火通信", "海格通信", "顺灏股份"),
    @("神剑股份", "神剑股份", "通宇通讯"),
    @("平潭发展", "五洲新春", "创新医疗"),
    @("航天电子", "平潭发展", "山子高科")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
